# Auto-generated edit script applying numeric corrections to Sheets
# per commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 996.75
$ws.Range("I41").Value = 996.75
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 996.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -556.75
$ws.Range("H87").Value = 91977.5
$ws.Range("J87").Value = 91977.5
$ws.Range("L87").Value = 91977.5
$ws.Range("N87").Value = -94473.5
$ws.Range("H90").Value = 91977.5
$ws.Range("J90").Value = 91977.5
$ws.Range("L90").Value = 275932.5
$ws.Range("N90").Value = -288412.5
$ws.Range("H94").Value = 1324.1666
$ws.Range("I94").Value = 1324.1666
$ws.Range("K94").Value = 1324.1666
$ws.Range("M94").Value = -873.1666
$ws.Range("H106").Value = 4707
$ws.Range("I106").Value = 2985
$ws.Range("K106").Value = 2985
$ws.Range("M106").Value = -2354
$ws.Range("H129").Value = 1243.2727
$ws.Range("I129").Value = 1020.6667
$ws.Range("J129").Value = 2245
$ws.Range("K129").Value = 3062.0001
$ws.Range("L129").Value = 6735
$ws.Range("M129").Value = 1937.9999
$ws.Range("N129").Value = -16735
$ws.Range("H137").Value = 5895.1177
$ws.Range("I137").Value = 5570.654
$ws.Range("K137").Value = 16711.962
$ws.Range("M137").Value = -14161.962
$ws.Range("H138").Value = 2540.5688
$ws.Range("J138").Value = 2979.6223
$ws.Range("L138").Value = 8938.866900000001
$ws.Range("N138").Value = -19218.8669
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21668434
$ws.Range("I74").Value = 41668332
$ws.Range("J74").Value = 1668538.9
$ws.Range("K74").Value = 41668332
$ws.Range("L74").Value = 1668538.9
$ws.Range("M74").Value = -41667458
$ws.Range("N74").Value = -1670286.9
$ws.Range("H77").Value = 21668434
$ws.Range("I77").Value = 41668332
$ws.Range("J77").Value = 1668538.9
$ws.Range("K77").Value = 208341660
$ws.Range("L77").Value = 8342694.5
$ws.Range("M77").Value = -208337292
$ws.Range("N77").Value = -8351430.5
$ws.Range("H102").Value = 14083.823
$ws.Range("I102").Value = 14798.25
$ws.Range("J102").Value = 12369.2
$ws.Range("K102").Value = 14798.25
$ws.Range("L102").Value = 12369.2
$ws.Range("M102").Value = -13176.25
$ws.Range("N102").Value = -15613.2
$ws.Range("H122").Value = 3670.7407
$ws.Range("I122").Value = 2643.7856
$ws.Range("J122").Value = 4776.6924
$ws.Range("K122").Value = 7931.3568
$ws.Range("L122").Value = 14330.0772
$ws.Range("M122").Value = -5481.3568
$ws.Range("N122").Value = -19230.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 865514.9
$ws.Range("I31").Value = 2523.3333
$ws.Range("J31").Value = 2591498
$ws.Range("K31").Value = 2523.3333
$ws.Range("L31").Value = 2591498
$ws.Range("M31").Value = -2228.3333
$ws.Range("N31").Value = -2592088
$ws.Range("H34").Value = 865514.9
$ws.Range("I34").Value = 2523.3333
$ws.Range("J34").Value = 2591498
$ws.Range("K34").Value = 2523.3333
$ws.Range("L34").Value = 2591498
$ws.Range("M34").Value = -2321.3333
$ws.Range("N34").Value = -2591902
$ws.Range("H92").Value = 49495
$ws.Range("J92").Value = 49495
$ws.Range("L92").Value = 49495
$ws.Range("N92").Value = -54487
$ws.Range("H107").Value = 2650.6667
$ws.Range("I107").Value = 1323.375
$ws.Range("K107").Value = 1323.375
$ws.Range("M107").Value = 596.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 807.3
$ws.Range("I23").Value = 1335.25
$ws.Range("J23").Value = 455.33334
$ws.Range("K23").Value = 4005.75
$ws.Range("L23").Value = 1366.00002
$ws.Range("M23").Value = -3770.75
$ws.Range("N23").Value = -1836.00002
$ws.Range("H37").Value = 67249.5
$ws.Range("J37").Value = 67249.5
$ws.Range("L37").Value = 201748.5
$ws.Range("N37").Value = -201972.5
$ws.Range("H51").Value = 40005
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H107").Value = 1634.25
$ws.Range("J107").Value = 1634.25
$ws.Range("L107").Value = 4902.75
$ws.Range("N107").Value = -8742.75
$ws.Range("H131").Value = 7261.83
$ws.Range("J131").Value = 7261.83
$ws.Range("L131").Value = 21785.49
$ws.Range("N131").Value = -31865.49
$ws.Range("H137").Value = 4787.357
$ws.Range("I137").Value = 4118.2856
$ws.Range("J137").Value = 5456.4287
$ws.Range("K137").Value = 12354.8568
$ws.Range("L137").Value = 16369.2861
$ws.Range("M137").Value = -7254.856800000001
$ws.Range("N137").Value = -26569.2861
$ws.Range("H140").Value = 7018.5
$ws.Range("I140").Value = 6429.6
$ws.Range("K140").Value = 19288.8
$ws.Range("M140").Value = -14108.8
$ws.Range("M51").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2912.2
$ws.Range("I80").Value = 2412.55
$ws.Range("K80").Value = 2412.55
$ws.Range("M80").Value = -1414.55
$ws.Range("H83").Value = 2912.2
$ws.Range("I83").Value = 2412.55
$ws.Range("K83").Value = 12062.75
$ws.Range("M83").Value = -7070.75
$ws.Range("H102").Value = 3270.2
$ws.Range("I102").Value = 2765.2727
$ws.Range("K102").Value = 2765.2727
$ws.Range("M102").Value = -1143.2727
$ws.Range("H114").Value = 31330
$ws.Range("J114").Value = 31330
$ws.Range("L114").Value = 31330
$ws.Range("N114").Value = -40008
$ws.Range("H126").Value = 4166.7617
$ws.Range("I126").Value = 3550.1667
$ws.Range("K126").Value = 10650.5001
$ws.Range("M126").Value = -8180.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 56488.9
$ws.Range("J7").Value = 91257.586
$ws.Range("L7").Value = 91257.586
$ws.Range("N7").Value = -91481.586
$ws.Range("H82").Value = 1330.9286
$ws.Range("I82").Value = 1270.4445
$ws.Range("J82").Value = 1439.8
$ws.Range("K82").Value = 1270.4445
$ws.Range("L82").Value = 1439.8
$ws.Range("M82").Value = -909.4445000000001
$ws.Range("N82").Value = -2161.8
$ws.Range("H85").Value = 1330.9286
$ws.Range("I85").Value = 1270.4445
$ws.Range("J85").Value = 1439.8
$ws.Range("K85").Value = 1270.4445
$ws.Range("L85").Value = 1439.8
$ws.Range("M85").Value = -22.44450000000006
$ws.Range("N85").Value = -3935.8
$ws.Range("H122").Value = 6985.0557
$ws.Range("I122").Value = 6295.615
$ws.Range("J122").Value = 8777.6
$ws.Range("K122").Value = 18886.845
$ws.Range("L122").Value = 26332.8
$ws.Range("M122").Value = -16436.845
$ws.Range("N122").Value = -31232.8
$ws.Range("H126").Value = 56488.9
$ws.Range("J126").Value = 91257.586
$ws.Range("L126").Value = 273772.758
$ws.Range("N126").Value = -278712.758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2792.8928
$ws.Range("I122").Value = 2766.7083
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 8300.124899999999
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -5850.124899999999
$ws.Range("N122").Value = -13750
